# Add "hydrogen combined cycle" as a power plant type on the
# RAF-generation sheet, renaming the existing "hydrogen" row to
# "hydrogen combustion turbine" in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAF-generation")

# Rename the existing "hydrogen" entry (row 24) to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add a new row for "hydrogen combined cycle" with the same 0.9 availability factor
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 0.9

# Reflect the new used range on the sheet and select the two hydrogen rows,
# matching the authored selection after the edit
$ws.Range("A24:A25").Select()
